# 16th may 2025- updated
# Adds four more product/price rows under the existing "Swag Labs" /
# "Sauce Labs Backpack" rows on Sheet2: price, price, qty, qty (row 6 is
# left blank, matching the authored edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# Leading apostrophe -> stored as text (quote-prefixed), matching the
# shared-string / quotePrefix cells in the target sheet. Assignment order
# controls the shared-string table order (29.99, 129.94, 6, 55).
$ws.Range("A4").Value = "'29.99"
$ws.Range("A5").Value = "'129.94"
$ws.Range("A3").Value = "'6"
$ws.Range("A7").Value = "'55"

# A4 ("29.99") carries its own font (pasted-in formatting): Arial 10pt,
# RGB 13,23,34 (hex FF132322).
$ws.Range("A4").Font.Name = "Arial"
$ws.Range("A4").Font.Size = 10
$ws.Range("A4").Font.Color = 2237203

# Column A widens from the old best-fit 18.09 to a fixed 21 characters.
$ws.Range("A1").ColumnWidth = 20.14

# Selection moves to C8.
[void]$ws.Range("C8").Select()

# Page orientation explicitly set to portrait.
$ws.PageSetup.Orientation = 1
